$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'320.10"
$ws.Range("E2").Formula = "'4.89%"

$ws.Range("D3").Formula = "'49.45"
$ws.Range("E3").Formula = "'11.57%"

$ws.Range("D4").Formula = "'5.322"
$ws.Range("E4").Formula = "'3.99%"

$ws.Range("D5").Formula = "'0.08059"
$ws.Range("E5").Formula = "'3.21%"

$ws.Range("D6").Formula = "'4.607"
$ws.Range("E6").Formula = "'4.31%"

$ws.Range("D7").Formula = "'1.344"
$ws.Range("E7").Formula = "'28.02%"

$ws.Range("D8").Formula = "'1.653"
$ws.Range("E8").Formula = "'1.98%"

$ws.Range("D9").Formula = "'0.1283"
$ws.Range("E9").Formula = "'-1.35%"

$ws.Range("D10").Formula = "'0.1961"
$ws.Range("E10").Formula = "'5.19%"

$ws.Range("D11").Formula = "'0.09564"
$ws.Range("E11").Formula = "'3.94%"

$ws.Range("D12").Formula = "'0.04533"
$ws.Range("E12").Formula = "'9.43%"

$ws.Range("D13").Formula = "'0.1048"
$ws.Range("E13").Formula = "'0.26%"

$ws.Range("D14").Formula = "'0.001305"
$ws.Range("E14").Formula = "'1.95%"

$ws.Range("D15").Formula = "'0.04212"
$ws.Range("E15").Formula = "'1.11%"

$ws.Range("D16").Formula = "'0.005924"
$ws.Range("E16").Formula = "'2.54%"

$ws.Range("E17").Formula = "'-0.15%"

$ws.Range("D18").Formula = "'2.473"
$ws.Range("E18").Formula = "'5.51%"

$ws.Range("E19").Formula = "'4.31%"

$ws.Range("D20").Formula = "'8.208"
$ws.Range("E20").Formula = "'2.44%"

$ws.Range("D21").Formula = "'0.1391"
$ws.Range("E21").Formula = "'2.11%"

$ws.Range("E22").Formula = "'9.69%"

$ws.Range("D23").Formula = "'0.001295"
$ws.Range("E23").Formula = "'1.53%"

$ws.Range("D24").Formula = "'0.004231"
$ws.Range("E24").Formula = "'-4.67%"

$ws.Range("D25").Formula = "'0.0001353"
$ws.Range("E25").Formula = "'0.77%"

$ws.Range("D26").Formula = "'0.0003545"

$ws.Range("D38").Formula = "'0.02696"
$ws.Range("E38").Formula = "'6.33%"

$ws.Range("D39").Formula = "'0.05931"
$ws.Range("E39").Formula = "'11.07%"

$ws.Range("D40").Formula = "'0.01082"
$ws.Range("E40").Formula = "'92.29%"

$ws.Range("D41").Formula = "'0.008048"
$ws.Range("E41").Formula = "'4.53%"

$ws.Range("D42").Formula = "'0.1465"
$ws.Range("E42").Formula = "'7.07%"

$ws.Range("D43").Formula = "'0.007535"
$ws.Range("E43").Formula = "'2.73%"

$ws.Range("D44").Formula = "'0.007935"
$ws.Range("E44").Formula = "'-4.66%"

$ws.Range("D45").Formula = "'0.3213"
$ws.Range("E45").Formula = "'6.48%"

$ws.Range("D46").Formula = "'0.00007042"
$ws.Range("E46").Formula = "'5.38%"

$ws.Range("D47").Formula = "'0.00000000751"
$ws.Range("E47").Formula = "'0.74%"

$ws.Range("D48").Formula = "'0.05594"
$ws.Range("E48").Formula = "'-7.38%"

$ws.Range("D49").Formula = "'0.004007"
$ws.Range("E49").Formula = "'0.74%"

$ws.Range("D50").Formula = "'0.00002104"
$ws.Range("E50").Formula = "'0.74%"

$ws.Range("D51").Formula = "'0.0002003"
$ws.Range("E51").Formula = "'0.74%"
